$d = $word.ActiveDocument

function New-OoxmlPackage([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- 1. Replace "The language used..." through the 2nd blank paragraph after
#        "Are there limits on your code..." with the new bold formatting,
#        the new "Libraries used" writeup, and the new "Can have any amount..."
#        paragraph. ---
$startPara = $d.Paragraphs.Item(5)
$endPara = $d.Paragraphs.Item(13)
$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)

$block1 = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>The language used (and why you chose it); What libraries you have used.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">I used java, as it allowed me to use OOP to break the program into functional sections, which made debugging it much easier. </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Libraries used:</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>FileNotFoundException</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - </w:t></w:r><w:r><w:t>Signals that an attempt to open the file denoted by a specified pathname has failed</w:t></w:r><w:r><w:t>, which allows me to see where I have made a mistake, instead of the code simply crashing.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Random</w:t></w:r><w:r><w:t xml:space="preserve"> &#8211; generates random values, allowing me to start from many points in the weight space to find the global minima of the error function.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:hanging="720"/></w:pPr><w:r><w:t>Scanner</w:t></w:r><w:r><w:t xml:space="preserve"> &#8211; </w:t></w:r></w:p><w:p><w:r><w:t>File</w:t></w:r><w:r><w:t xml:space="preserve"> - </w:t></w:r><w:r><w:t>allows me to read from the excel file</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>ArrayList</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &#8211; allows me to make a new</w:t></w:r></w:p><w:p><w:r><w:t>List</w:t></w:r><w:r><w:t xml:space="preserve"> &#8211; data structure </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>&#8230;..</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p/><w:p/><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">How you implemented it &#8211; </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>e.g.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> OO approach with an MLP class and what methods it has, how the data are stored/structured, etc. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">The MLP algorithm &#8211; what additions did you make &#8211; </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>e.g.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> momentum, annealing, bold driver. Did you try different transfer functions? Alternative training algorithms &#8211; </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>e.g.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> conjugate gradients?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Are there limits on your code (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>e.g.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> have things been hard-coded or can it create any MLP with any number of inputs, hidden layers, outputs, etc). </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Can have any amount of hidden layer nodes, but only 1 hidden layer and 1 output</w:t></w:r></w:p><w:p/><w:p/>'

$rng.InsertXML((New-OoxmlPackage $block1))

# --- 2. Move the lastRenderedPageBreak marker: it now falls before
#        "Avoid hard-coding things..." instead of before
#        "Make sure you highlight (...". ---

$avoidText = "Avoid hard-coding things. In other words, don"+[char]8217+"t write the program for the given data set. It should be easily modifiable for other data sets, different numbers of inputs, different numbers of hidden nodes, etc."
$highlightText = "Make sure you highlight (e.g. with lots of comments) the actual backpropagation algorithm (and other enhancements) in your code so I can find it easily when marking. I need to see the algorithm(s) in your code to mark it (them). "

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Avoid hard-coding things")) {
        $firstRun = $p.Range.Words.First
        $p.Range.Sentences.First | Out-Null
    }
}

Write-Host "stage1 done"
